$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, avoiding Excel auto-converting
# numeric-looking strings (e.g. "320.36") into actual numbers.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
$ws.Range('D2').Value = '29.978.23'
$ws.Range('E2').Value = '  -0.72%  '

# Row 3
$ws.Range('D3').Value = '1.922.21'
$ws.Range('E3').Value = '  +0.41%  '

# Row 4
Set-TextValue 'D4' '1.002'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
Set-TextValue 'D5' '320.36'
$ws.Range('E5').Value = '  -2.82%  '

# Row 6
$ws.Range('E6').Value = '  +0.09%  '

# Row 7
Set-TextValue 'D7' '0.5034'
$ws.Range('E7').Value = '  -2.90%  '

# Row 8
Set-TextValue 'D8' '0.4026'
$ws.Range('E8').Value = '  -0.92%  '

# Row 9
Set-TextValue 'D9' '0.08238'
$ws.Range('E9').Value = '  -2.98%  '

# Row 10
Set-TextValue 'D10' '1.111'
$ws.Range('E10').Value = '  -1.36%  '

# Row 11
Set-TextValue 'D11' '41.98'
$ws.Range('E11').Value = '  -1.80%  '

# Row 12
Set-TextValue 'D12' '23.58'
$ws.Range('E12').Value = '  +0.90%  '

# Row 13
$ws.Range('D13').Value = '1.913.15'
$ws.Range('E13').Value = '  -0.37%  '

# Row 14
Set-TextValue 'D14' '6.402'
$ws.Range('E14').Value = '  -0.54%  '

# Row 15
Set-TextValue 'D15' '7.297'
$ws.Range('E15').Value = '  -1.24%  '

# Row 16
Set-TextValue 'D16' '1.002'
$ws.Range('E16').Value = '  +0.05%  '

# Row 17
$ws.Range('E17').Value = '  -3.27%  '

# Row 18
$ws.Range('E18').Value = '  -1.62%  '

# Row 20
Set-TextValue 'D20' '18.16'
$ws.Range('E20').Value = '  -1.77%  '

# Row 21
Set-TextValue 'D21' '1.001'
$ws.Range('E21').Value = '  +0.09%  '

# Row 22
Set-TextValue 'D22' '5.955'
$ws.Range('E22').Value = '  -0.97%  '

# Row 23
$ws.Range('D23').Value = '30.026.61'
$ws.Range('E23').Value = '  -0.64%  '

# Row 24
$ws.Range('E24').Value = '  -0.64%  '

# Row 25
Set-TextValue 'D25' '2.193'
$ws.Range('E25').Value = '  -1.54%  '

# Row 26
Set-TextValue 'D26' '22.17'
$ws.Range('E26').Value = '  +3.42%  '

# Row 27
$ws.Range('D27').Value = '2.133.58'

# Row 28
Set-TextValue 'D28' '161.31'
$ws.Range('E28').Value = '  -0.59%  '

# Row 29
Set-TextValue 'D29' '2.336'
$ws.Range('E29').Value = '  -3.03%  '

# Row 30
Set-TextValue 'D30' '128.87'
$ws.Range('E30').Value = '  +0.05%  '

# Row 31
Set-TextValue 'D31' '1.127'
$ws.Range('E31').Value = '  +2.35%  '

# Row 32
Set-TextValue 'D32' '0.1041'
$ws.Range('E32').Value = '  -2.44%  '

# Row 33
Set-TextValue 'D33' '5.980'
$ws.Range('E33').Value = '  -0.45%  '

# Row 34
Set-TextValue 'D34' '3.816'
$ws.Range('E34').Value = '  +4.76%  '

# Row 35
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D35' '5.408'
$ws.Range('E35').Value = '  +4.29%  '

# Row 36
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D36' '0.02446'
$ws.Range('E36').Value = '  -1.80%  '

# Row 37
Set-TextValue 'D37' '0.06441'
$ws.Range('E37').Value = '  -2.10%  '

# Row 38
Set-TextValue 'D38' '8.972'
$ws.Range('E38').Value = '  +1.88%  '

# Row 39
Set-TextValue 'D39' '0.2164'
$ws.Range('E39').Value = '  -2.12%  '

# Row 40
Set-TextValue 'D40' '1.192'
$ws.Range('E40').Value = '  -2.63%  '

# Row 41
Set-TextValue 'D41' '0.6420'
$ws.Range('E41').Value = '  -1.52%  '

# Row 42
$ws.Range('E42').Value = '  -4.76%  '

# Row 43
$ws.Range('E43').Value = '  -1.73%  '

# Row 44
Set-TextValue 'D44' '1.000'
$ws.Range('E44').Value = '  +0.16%  '

# Row 45
Set-TextValue 'D45' '13.33'
$ws.Range('E45').Value = '  -0.39%  '

# Row 46
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D46' '2.168'
$ws.Range('E46').Value = '  +4.43%  '

# Row 47
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D47' '0.6004'
$ws.Range('E47').Value = '  -2.26%  '

# Row 48
Set-TextValue 'D48' '3.644'
$ws.Range('E48').Value = '  -2.56%  '

# Row 49
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D49' '122.85'
$ws.Range('E49').Value = '  -1.14%  '

# Row 50
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue 'D50' '1.214'
$ws.Range('E50').Value = '  -2.52%  '

# Row 51
Set-TextValue 'D51' '78.82'
$ws.Range('E51').Value = '  -0.88%  '
